$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new "aula 11" entry ---------------------------------------
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = "2. Inciando o Desenvolvimento do Projeto"
$ws.Range("D6").Value = "11. Incluíndo o recurso de WebJars"
$ws.Range("E6").Value = '6:08 - foi ensinado uma forma de adicionar bibliotecas para pagina HTML (Jquery, bootstrap, icones) através de dependências no pom.xml. Usando esta forma, o "src" do documento html deve direcionar para o diretorio raiz onde encontra-se os webjars baixados pelo pom. É possivel ver o local destes diretórios direto na documentação das bibliotecas. Os webjars podem ser baixados atraves do site https://www.webjars.org/'

# Same wrap-text look as the other "observação" cells (E4/E5) and a taller row
$ws.Range("E6").WrapText = $true
$ws.Rows(6).RowHeight = 90

# --- Row 7: empty cell that just carries an underlined-font style -----
$ws.Range("E7").Font.Underline = $true

# --- Selection left on E7, matching the saved workbook ----------------
$ws.Range("E7").Select() | Out-Null

# --- Page setup (paper size / orientation) -----------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
